$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows 10-16 (existing rows, content reshuffled + new values) ---
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 0.9859038466344734
$ws.Range("D10").Value = 1.017880101365173
$ws.Range("E10").Value = 0.9888596471693042
$ws.Range("F10").Value = 0.9859038466344734
$ws.Range("G10").Value = 1.008540572174369
$ws.Range("H10").Value = 0.981394061939606
$ws.Range("I10").Value = 0.9882060459891621
$ws.Range("J10").Value = 1.017880101365173
$ws.Range("K10").Value = 1.003369874267239
$ws.Range("L10").Value = 0.994636860450856
$ws.Range("M10").Value = 0.995130712545348

$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.9977137665264835
$ws.Range("D11").Value = 0.9508358715485322
$ws.Range("E11").Value = 1.005980228825406
$ws.Range("F11").Value = 0.9977137665264835
$ws.Range("G11").Value = 0.9698483517896725
$ws.Range("H11").Value = 1.029504014277771
$ws.Range("I11").Value = 1.004760138678962
$ws.Range("J11").Value = 0.9508358715485322
$ws.Range("K11").Value = 0.9784080501869692
$ws.Range("L11").Value = 0.9880609083567262
$ws.Range("M11").Value = 0.9931070619411381

$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.9976010191376431
$ws.Range("D12").Value = 0.9513918949420904
$ws.Range("E12").Value = 1.005829722068599
$ws.Range("F12").Value = 0.9976010191376431
$ws.Range("G12").Value = 0.9701422508012693
$ws.Range("H12").Value = 1.029158336382481
$ws.Range("I12").Value = 1.004626759878114
$ws.Range("J12").Value = 0.9513918949420904
$ws.Range("K12").Value = 0.9786108085053444
$ws.Range("L12").Value = 0.9881059138214938
$ws.Range("M12").Value = 0.9931249972016994

$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.997665206992058
$ws.Range("D13").Value = 0.9509960216976289
$ws.Range("E13").Value = 1.005916979274208
$ws.Range("F13").Value = 0.997665206992058
$ws.Range("G13").Value = 0.9699099287467844
$ws.Range("H13").Value = 1.029487495263738
$ws.Range("I13").Value = 1.004720465747753
$ws.Range("J13").Value = 0.9509960216976289
$ws.Range("K13").Value = 0.9784565004859185
$ws.Range("L13").Value = 0.9880608537389881
$ws.Range("M13").Value = 0.9931160162870284

$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.9735559999999986
$ws.Range("D14").Value = 1.055587999999999
$ws.Range("E14").Value = 0.9833439999999984
$ws.Range("F14").Value = 0.9735559999999986
$ws.Range("G14").Value = 1.032276
$ws.Range("H14").Value = 0.9423
$ws.Range("I14").Value = 0.9768839999999996
$ws.Range("J14").Value = 1.055587999999999
$ws.Range("K14").Value = 1.019465999999999
$ws.Range("L14").Value = 0.9965109999999987
$ws.Range("M14").Value = 0.9939913333333328

$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0.97
$ws.Range("D15").Value = 1.08
$ws.Range("E15").Value = 0.9773874999999994
$ws.Range("F15").Value = 0.97
$ws.Range("G15").Value = 1.045399999999998
$ws.Range("H15").Value = 0.91
$ws.Range("I15").Value = 0.97
$ws.Range("J15").Value = 1.08
$ws.Range("K15").Value = 1.02869375
$ws.Range("L15").Value = 0.9993468749999999
$ws.Range("M15").Value = 0.9921312499999996

$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.9801732665344044
$ws.Range("D16").Value = 1.044343337574397
$ws.Range("E16").Value = 0.9842803031040035
$ws.Range("F16").Value = 0.9801732665344044
$ws.Range("G16").Value = 1.023676689510398
$ws.Range("H16").Value = 0.9470913673216019
$ws.Range("I16").Value = 0.9802400616448033
$ws.Range("J16").Value = 1.044343337574397
$ws.Range("K16").Value = 1.0143118203392
$ws.Range("L16").Value = 0.9972425434368025
$ws.Range("M16").Value = 0.9933008376149348

# --- Add new rows 17-19 ---
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9945762888340077
$ws.Range("D17").Value = 0.9942774305556739
$ws.Range("E17").Value = 0.9945914620329442
$ws.Range("F17").Value = 0.9945762888340077
$ws.Range("G17").Value = 0.9942513311365048
$ws.Range("H17").Value = 0.9943527493820711
$ws.Range("I17").Value = 0.9940393577652945
$ws.Range("J17").Value = 0.9942774305556739
$ws.Range("K17").Value = 0.9944344462943091
$ws.Range("L17").Value = 0.9945053675641584
$ws.Range("M17").Value = 0.994348103284416

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 0.9951608369074652
$ws.Range("D18").Value = 0.9909266503011143
$ws.Range("E18").Value = 0.994926238835801
$ws.Range("F18").Value = 0.9951608369074652
$ws.Range("G18").Value = 0.9929793878313541
$ws.Range("H18").Value = 0.9954265536492175
$ws.Range("I18").Value = 0.9949208831756374
$ws.Range("J18").Value = 0.9909266503011143
$ws.Range("K18").Value = 0.9929264445684576
$ws.Range("L18").Value = 0.9940436407379614
$ws.Range("M18").Value = 0.9940567584500982

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.9966544266109879
$ws.Range("D19").Value = 0.9842904832345777
$ws.Range("E19").Value = 0.9969959913512131
$ws.Range("F19").Value = 0.9966544266109879
$ws.Range("G19").Value = 0.9883209725765859
$ws.Range("H19").Value = 1.001258586808111
$ws.Range("I19").Value = 0.9968032194135852
$ws.Range("J19").Value = 0.9842904832345777
$ws.Range("K19").Value = 0.9906432372928954
$ws.Range("L19").Value = 0.9936488319519416
$ws.Range("M19").Value = 0.9940539466658435

# Copy formatting (bold/border/center) from row 16 col A to new rows 17-19 col A
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17:A19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

Write-Host "Done"